# Generate Report for Archive
# -----------------------------------------------------------------------
# The localization report was regenerated and two entries
#   05294cb7-e2f2-411b-a2bd-ca4347d00657.md   (status: Ready for handoff)
#   fc368083-54a6-4157-a90d-0fc2c1a1fe53.md   (status: In Translation)
# swapped places (fc368083 now sorts above 05294cb7) on every sheet
# (Overview, zh-cn, de-de). This script rewrites the affected cells and
# re-creates the hyperlink list for each sheet to match the regenerated
# report, row by row, top to bottom.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 4 now holds fc368083's entry, which is "In Translation".
$ws1.Range("A4").Value = "fc368083-54a6-4157-a90d-0fc2c1a1fe53.md"
$ws1.Range("B4").Value = "In Translation"
$ws1.Range("C4").Value = "In Translation"

# Row 5 now holds 05294cb7's entry, which is "Ready for handoff".
$ws1.Range("A5").Value = "05294cb7-e2f2-411b-a2bd-ca4347d00657.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

# Rebuild the hyperlink list top to bottom so relationship ids stay in
# the same sequential order the report generator produces.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/dbb8b7c2-6efa-4117-90da-56923094cd06.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9c78ddd6f28012dea28437b984c3a8b47508a7b3/e2e/05294cb7-e2f2-411b-a2bd-ca4347d00657.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f4055970680395dc3b7eecf6d9ddfc37e7948de6/e2e/fc368083-54a6-4157-a90d-0fc2c1a1fe53.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/9c78ddd6f28012dea28437b984c3a8b47508a7b3/.localization-config") | Out-Null

# ---- zh-cn sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A4").Value = "fc368083-54a6-4157-a90d-0fc2c1a1fe53.md"
$ws2.Range("B4").Value = "In Translation"
$ws2.Range("C4").Value = "fc368083-54a6-4157-a90d-0fc2c1a1fe53.17ffdc14645d6678fc3687045001f2fb09a8a6dc.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-07 09:30:28"

$ws2.Range("A5").Value = "05294cb7-e2f2-411b-a2bd-ca4347d00657.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = "05294cb7-e2f2-411b-a2bd-ca4347d00657.68f81ef12b8a67246ab4a8ce9182b2002eeb404d.zh-cn.xlf"
$ws2.Range("D5").Value = "2016-03-07 09:31:20"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f566921d63ecb759a7405dc81c59ef598af8399/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.0e70ba733ba7709f7ac6be56e782aba445bc8037.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/dbb8b7c2-6efa-4117-90da-56923094cd06.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f566921d63ecb759a7405dc81c59ef598af8399/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/dbb8b7c2-6efa-4117-90da-56923094cd06.9101628fb3d38b6c45007db9c130d59806c44cb6.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9c78ddd6f28012dea28437b984c3a8b47508a7b3/e2e/05294cb7-e2f2-411b-a2bd-ca4347d00657.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/28aeefee9c0399befa5c7cf511ee7493db326010/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/05294cb7-e2f2-411b-a2bd-ca4347d00657.68f81ef12b8a67246ab4a8ce9182b2002eeb404d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f4055970680395dc3b7eecf6d9ddfc37e7948de6/e2e/fc368083-54a6-4157-a90d-0fc2c1a1fe53.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c21722051eaa0afb7f3f87cf3650a9a99b61a757/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fc368083-54a6-4157-a90d-0fc2c1a1fe53.17ffdc14645d6678fc3687045001f2fb09a8a6dc.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/9c78ddd6f28012dea28437b984c3a8b47508a7b3/.localization-config") | Out-Null

# ---- de-de sheet -------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A4").Value = "fc368083-54a6-4157-a90d-0fc2c1a1fe53.md"
$ws3.Range("B4").Value = "In Translation"
$ws3.Range("C4").Value = "fc368083-54a6-4157-a90d-0fc2c1a1fe53.17ffdc14645d6678fc3687045001f2fb09a8a6dc.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-07 09:30:40"

$ws3.Range("A5").Value = "05294cb7-e2f2-411b-a2bd-ca4347d00657.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = "05294cb7-e2f2-411b-a2bd-ca4347d00657.68f81ef12b8a67246ab4a8ce9182b2002eeb404d.de-de.xlf"
$ws3.Range("D5").Value = "2016-03-07 09:31:34"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b2e0e8513621865dd16d434a18b4bcf509d4fbc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.0e70ba733ba7709f7ac6be56e782aba445bc8037.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/dbb8b7c2-6efa-4117-90da-56923094cd06.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b2e0e8513621865dd16d434a18b4bcf509d4fbc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/dbb8b7c2-6efa-4117-90da-56923094cd06.9101628fb3d38b6c45007db9c130d59806c44cb6.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9c78ddd6f28012dea28437b984c3a8b47508a7b3/e2e/05294cb7-e2f2-411b-a2bd-ca4347d00657.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/413d8b59fa90258f9243ba0af727baca699eb31c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/05294cb7-e2f2-411b-a2bd-ca4347d00657.68f81ef12b8a67246ab4a8ce9182b2002eeb404d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f4055970680395dc3b7eecf6d9ddfc37e7948de6/e2e/fc368083-54a6-4157-a90d-0fc2c1a1fe53.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b2db43b3d259b99f5ee21ca15eecddda9566fe4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fc368083-54a6-4157-a90d-0fc2c1a1fe53.17ffdc14645d6678fc3687045001f2fb09a8a6dc.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/9c78ddd6f28012dea28437b984c3a8b47508a7b3/.localization-config") | Out-Null
